# Build site at 2023-04-12 14:53:07 UTC
# Applies the LOQ4263.xlsx content update:
#  - Inserts a new row (the "Docentes responsaveis:" / professor name row) at row 13,
#    pushing the existing rows 13-23 down to 14-24.
#  - Fills in the Portuguese translations / texts that were previously missing
#    for Objetivos, Programa resumido, Programa, Metodo, Criterio,
#    Norma de recuperacao and Bibliografia.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row at 13 - this shifts the old rows 13..23 down to 14..24
#    and carries the row heights / cell styles along automatically (matches
#    the row-height pattern seen in the target worksheet).
$ws.Rows(13).Insert()

# 2) The new row 13 becomes the "Docentes responsaveis:" value row - it shows
#    the professor's name in columns B and C (column A keeps the
#    "Docentes responsaveis:" label that already lived in A12).
$ws.Range("B13").Value = "5840535 - Messias Borges Silva"
$ws.Range("C13").Value = "5840535 - Messias Borges Silva"

# 3) Row 10 (Objetivos:) - fill in the Portuguese objectives text that was missing.
$objetivosPt = "Explicar os conceitos, métodos e resolver problemas que ilustrem aplicações sem recorrer a desenvolvimento teóricos da Teoria da Confiabilidade. Pretende-se uma formação geral com o uso de modelos probabilísticos e estatísticos, e com aplicações na área de engenharia. Uso de aplicativos computacionais para análise de conjunto de dados."
$ws.Range("B10").Value = $objetivosPt
$ws.Range("C10").Value = $objetivosPt

# 4) Row 14 (Programa resumido:) - fill in the Portuguese short-syllabus text.
$shortSyllabusPt = "1. Confiabilidade e disponibilidade de sistemas. 2. Famílias de distribuições. 3. Sistemas reparáveis. 4. Análise gráfica de dados. 5. Estimação de características de um sistema. 6. FMEA, 7. Aplicações na Gestão da Manutenção, 8. Manutenção Produtiva Total, 9. Design for Six Sigma. 10. RCM (Reliability Centered Maintenance)"
$ws.Range("B14").Value = $shortSyllabusPt
$ws.Range("C14").Value = $shortSyllabusPt

# 5) Row 16 (Programa:) - fill in the Portuguese full-syllabus text.
$syllabusPt = "1. Confiabilidade e disponibilidade de sistemas, decomposição por cortes e caminhos, árvores de eventos. 2. Famílias de distribuições úteis em Teoria a Confiabilidade. 3. Sistemas reparáveis, manutenção, aproximações assintóticas. 4. Análise gráfica de dados. 5. Estimação do tempo de vida e das características de um sistema. 6. FMEA, 7. Aplicações na Gestão da Manutenção, 8. Manutenção Produtiva Total, 9. Design for Six Sigma DFSS 10. RCM (Reliability Centered Maintenance)"
$ws.Range("B16").Value = $syllabusPt
$ws.Range("C16").Value = $syllabusPt

# 6) Row 19 (Metodo:) - fill in the evaluation method text.
$metodoPt = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."
$ws.Range("B19").Value = $metodoPt
$ws.Range("C19").Value = $metodoPt

# 7) Row 20 (Criterio:) - fill in the passing-grade criterion text.
$criterioPt = "NF≥ 5,0."
$ws.Range("B20").Value = $criterioPt
$ws.Range("C20").Value = $criterioPt

# 8) Row 21 (Norma de recuperacao:) - fill in the make-up exam rule text.
$normaPt = "Média aritmética da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recuperação."
$ws.Range("B21").Value = $normaPt
$ws.Range("C21").Value = $normaPt

# 9) Row 22 (Bibliografia:) - fill in the bibliography text.
$bibliografiaPt = "1. S. Zachs, INTRODUCTION TO RELIABILITY ANALYSIS: PROBABILITY MODELS AND STATISTICAL METHODS, Springer Verlag, New York, 19922. I.B. Gertsbakh, STATISTICAL RELIABILITY THEORY, Marcel Dekker, New York, 1989.3. J. Knezevic, RELIABILITY, MAINTAINABILITY, AND SUPPORTABILITY: A PROBABILITY APPROACH, McGraw-Hill, 1993.4. R.S. Dhillon, C. Singh, ENGINEERING RELIABILITY. NEW TECHNIQUES AND APPLICATIONS, Wiley Interscience, 1981. 5. HARRY, M. , LINSENMANND.R., The Six Sigma Fieldbook, Doubleday, New York, 2006"
$ws.Range("B22").Value = $bibliografiaPt
$ws.Range("C22").Value = $bibliografiaPt
